$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.47"
$ws.Range("E2").Value = "'-0.73%"
$ws.Range("E3").Value = "'-4.30%"
$ws.Range("D4").Value = "'5.262"
$ws.Range("E4").Value = "'1.61%"
$ws.Range("D5").Value = "'0.05702"
$ws.Range("E5").Value = "'-0.55%"
$ws.Range("D6").Value = "'6.639"
$ws.Range("E6").Value = "'0.98%"
$ws.Range("D7").Value = "'3.208"
$ws.Range("E7").Value = "'3.68%"
$ws.Range("D8").Value = "'0.8511"
$ws.Range("E8").Value = "'-0.77%"
$ws.Range("D9").Value = "'0.8864"
$ws.Range("E9").Value = "'1.93%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1391"
$ws.Range("E10").Value = "'1.95%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07084"
$ws.Range("E11").Value = "'0.28%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03140"
$ws.Range("E12").Value = "'7.83%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09228"
$ws.Range("E13").Value = "'-1.68%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001539"
$ws.Range("E14").Value = "'0.91%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0005989"
$ws.Range("E15").Value = "'-94.18%"
$ws.Range("D16").Value = "'0.005945"
$ws.Range("E16").Value = "'-0.28%"
$ws.Range("D17").Value = "'3.489"
$ws.Range("E17").Value = "'-0.02%"
$ws.Range("E18").Value = "'-0.07%"
$ws.Range("D19").Value = "'0.3167"
$ws.Range("E19").Value = "'-0.58%"
$ws.Range("D20").Value = "'0.03312"
$ws.Range("E20").Value = "'-2.30%"
$ws.Range("D21").Value = "'0.1294"
$ws.Range("E21").Value = "'0.62%"
$ws.Range("D22").Value = "'3.519"
$ws.Range("E22").Value = "'1.57%"
$ws.Range("D23").Value = "'0.04083"
$ws.Range("E23").Value = "'-1.00%"
$ws.Range("D24").Value = "'0.1379"
$ws.Range("D25").Value = "'0.001224"
$ws.Range("E25").Value = "'-0.15%"
$ws.Range("D26").Value = "'0.004151"
$ws.Range("E26").Value = "'-17.15%"
$ws.Range("D27").Value = "'0.0001200"
$ws.Range("D28").Value = "'0.0001449"
$ws.Range("D40").Value = "'0.03800"
$ws.Range("E40").Value = "'1.24%"
$ws.Range("D41").Value = "'0.1065"
$ws.Range("E41").Value = "'-0.69%"
$ws.Range("D42").Value = "'0.003739"
$ws.Range("E42").Value = "'7.40%"
$ws.Range("E43").Value = "'-10.19%"
$ws.Range("D44").Value = "'0.009489"
$ws.Range("E44").Value = "'12.04%"
$ws.Range("D45").Value = "'0.00005275"
$ws.Range("E45").Value = "'0.44%"
$ws.Range("E46").Value = "'0.02%"
$ws.Range("E47").Value = "'62.32%"
$ws.Range("E48").Value = "'-0.30%"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'0.02%"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'0.02%"
